$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels: rules.minimum -> rules.dateMinimum, rules.maximum -> rules.dateMaximum
$ws.Range("D1").Value = "rules.dateMinimum"
$ws.Range("E1").Value = "rules.dateMaximum"

# Update the selected/active cell from E8 to G9
$ws.Range("G9").Select()
